$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 93
$ws.Range("I9").Value = 0
$ws.Range("K9").Value = 0
$ws.Range("M9").ClearContents()
$ws.Range("H34").Value = 3511
$ws.Range("I34").Value = 1348
$ws.Range("J34").Value = 10000
$ws.Range("K34").Value = 1348
$ws.Range("L34").Value = 10000
$ws.Range("M34").Value = -1145
$ws.Range("N34").Value = -10406
$ws.Range("H36").Value = 3511
$ws.Range("I36").Value = 1348
$ws.Range("J36").Value = 10000
$ws.Range("K36").Value = 1348
$ws.Range("L36").Value = 10000
$ws.Range("M36").Value = -633
$ws.Range("N36").Value = -11430
$ws.Range("H76").Value = 3893.4092
$ws.Range("I76").Value = 3944.6875
$ws.Range("J76").Value = 3756.6667
$ws.Range("K76").Value = 3944.6875
$ws.Range("L76").Value = 3756.6667
$ws.Range("M76").Value = -3629.6875
$ws.Range("N76").Value = -4386.6667
$ws.Range("H79").Value = 3893.4092
$ws.Range("I79").Value = 3944.6875
$ws.Range("J79").Value = 3756.6667
$ws.Range("K79").Value = 3944.6875
$ws.Range("L79").Value = 3756.6667
$ws.Range("M79").Value = -2852.6875
$ws.Range("N79").Value = -5940.6667
$ws.Range("H107").Value = 7283.4287
$ws.Range("I107").Value = 8471.5
$ws.Range("J107").Value = 155
$ws.Range("K107").Value = 8471.5
$ws.Range("L107").Value = 155
$ws.Range("M107").Value = -6551.5
$ws.Range("N107").Value = -3995
$ws.Range("H137").Value = 7712348
$ws.Range("I137").Value = 890.82355
$ws.Range("J137").Value = 22278434
$ws.Range("K137").Value = 2672.47065
$ws.Range("L137").Value = 66835302
$ws.Range("M137").Value = -122.4706499999998
$ws.Range("N137").Value = -66840402
$ws.Range("H138").Value = 4018362.8
$ws.Range("I138").Value = 7753324.5
$ws.Range("J138").Value = 3278.925
$ws.Range("K138").Value = 23259973.5
$ws.Range("L138").Value = 9836.775000000001
$ws.Range("M138").Value = -23254833.5
$ws.Range("N138").Value = -20116.775
$ws.Range("H141").Value = 2653.8948
$ws.Range("I141").Value = 1571.0769
$ws.Range("K141").Value = 4713.2307
$ws.Range("M141").Value = 466.7692999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 797.53845
$ws.Range("I2").Value = 723.57574
$ws.Range("J2").Value = 1204.3334
$ws.Range("K2").Value = 723.57574
$ws.Range("L2").Value = 1204.3334
$ws.Range("M2").Value = -610.57574
$ws.Range("N2").Value = -1430.3334
$ws.Range("H61").Value = 11630045
$ws.Range("I61").Value = 13890960
$ws.Range("J61").Value = 2485.7144
$ws.Range("K61").Value = 13890960
$ws.Range("L61").Value = 2485.7144
$ws.Range("M61").Value = -13890748
$ws.Range("N61").Value = -2909.7144
$ws.Range("H116").Value = 797.53845
$ws.Range("I116").Value = 723.57574
$ws.Range("J116").Value = 1204.3334
$ws.Range("K116").Value = 723.57574
$ws.Range("L116").Value = 1204.3334
$ws.Range("M116").Value = 1570.42426
$ws.Range("N116").Value = -5792.3334
$ws.Range("H122").Value = 6768.04
$ws.Range("I122").Value = 8439.895
$ws.Range("J122").Value = 1473.8334
$ws.Range("K122").Value = 25319.685
$ws.Range("L122").Value = 4421.5002
$ws.Range("M122").Value = -22869.685
$ws.Range("N122").Value = -9321.5002
$ws.Range("H136").Value = 11630045
$ws.Range("I136").Value = 13890960
$ws.Range("J136").Value = 2485.7144
$ws.Range("K136").Value = 41672880
$ws.Range("L136").Value = 7457.1432
$ws.Range("M136").Value = -41670330
$ws.Range("N136").Value = -12557.1432

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 797.53845
$ws.Range("I3").Value = 723.57574
$ws.Range("J3").Value = 1204.3334
$ws.Range("K3").Value = 723.57574
$ws.Range("L3").Value = 1204.3334
$ws.Range("M3").Value = -609.57574
$ws.Range("N3").Value = -1432.3334
$ws.Range("H37").Value = 347.7143
$ws.Range("I37").Value = 389
$ws.Range("J37").Value = 100
$ws.Range("K37").Value = 389
$ws.Range("L37").Value = 100
$ws.Range("M37").Value = -252
$ws.Range("N37").Value = -374

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1131.1818
$ws.Range("J34").Value = 2294.6
$ws.Range("L34").Value = 6883.799999999999
$ws.Range("N34").Value = -7051.799999999999
$ws.Range("H138").Value = 6378.4614
$ws.Range("I138").Value = 1145.7142
$ws.Range("J138").Value = 12483.333
$ws.Range("K138").Value = 3437.1426
$ws.Range("L138").Value = 37449.999
$ws.Range("M138").Value = 1702.8574
$ws.Range("N138").Value = -47729.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 556388.4399999999
$ws.Range("I3").Value = 2500240
$ws.Range("J3").Value = 1002.2857
$ws.Range("K3").Value = 2500240
$ws.Range("L3").Value = 1002.2857
$ws.Range("M3").Value = -2500124
$ws.Range("N3").Value = -1234.2857
$ws.Range("H102").Value = 3781.6843
$ws.Range("I102").Value = 5237.6665
$ws.Range("K102").Value = 5237.6665
$ws.Range("M102").Value = -3615.6665
$ws.Range("H132").Value = 8405.105
$ws.Range("I132").Value = 10563.214
$ws.Range("J132").Value = 2362.4
$ws.Range("K132").Value = 31689.642
$ws.Range("L132").Value = 7087.200000000001
$ws.Range("M132").Value = -29159.642
$ws.Range("N132").Value = -12147.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3746.725
$ws.Range("I40").Value = 4348.0527
$ws.Range("K40").Value = 4348.0527
$ws.Range("M40").Value = -4212.0527
$ws.Range("H136").Value = 16672262
$ws.Range("I136").Value = 22729934
$ws.Range("J136").Value = 13663.75
$ws.Range("K136").Value = 68189802
$ws.Range("L136").Value = 40991.25
$ws.Range("M136").Value = -68187252
$ws.Range("N136").Value = -46091.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H43").Value = 12000
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 12000
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 12000
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -12298
$ws.Range("H136").Value = 818.30615
$ws.Range("I136").Value = 889.5641000000001
$ws.Range("J136").Value = 540.4
$ws.Range("K136").Value = 2668.6923
$ws.Range("L136").Value = 1621.2
$ws.Range("M136").Value = -118.6923000000002
$ws.Range("N136").Value = -6721.2
